# Apply updated NATMI Col1a2-Itgb1 LR-pair statistics ("Natmi following Dr Hou advice")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; F=1; G=10.45491533333333; H=31.364746;         I=0.0134573334963438; J=0.0134573334963438; K=3; L=1; M=117.044563;        N=351.133689;  O=0.3245365645427815; P=0.3245365645427815; Q=1223.690996391999;  R=11013.21896752799;  S=0.004367396780809915; T=0.004367396780809914 }
    3  = @{ E=3; F=1; G=10.45491533333333; H=31.364746;         I=0.0134573334963438; J=0.0134573334963438; K=3; L=1; M=101.5800373333333; N=304.740112;  O=0.281657135515876;  P=0.281657135515876;  Q=1062.010689876839;  R=9558.096208891553;   S=0.003790354004262044; T=0.003790354004262042 }
    4  = @{ E=3; F=1; G=10.45491533333333; H=31.364746;         I=0.0134573334963438; J=0.0134573334963438; K=3; L=1; M=142.0267893333333; N=426.080368;  O=0.3938062999413425; P=0.3938062999413425; Q=1484.87805754517;   R=13363.90251790653;  S=0.005299582711271842; T=0.005299582711271842 }
    5  = @{ E=3; F=1; G=735.4993083333334; H=2206.497925;       I=0.9467182815928301; J=0.9467182815928301; K=3; L=1; M=117.044563;        N=351.133689;  O=0.3245365645427815; P=0.3245365645427815; Q=86086.19513067725;  R=774775.7561760953;  S=0.3072446986979827;  T=0.3072446986979827  }
    6  = @{ E=3; F=1; G=735.4993083333334; H=2206.497925;       I=0.9467182815928301; J=0.9467182815928301; K=3; L=1; M=101.5800373333333; N=304.740112;  O=0.281657135515876;  P=0.281657135515876;  Q=74712.04719914086;  R=672408.4247922676;  S=0.266649959333949;   T=0.266649959333949   }
    7  = @{ E=3; F=1; G=735.4993083333334; H=2206.497925;       I=0.9467182815928301; J=0.9467182815928301; K=3; L=1; M=142.0267893333333; N=426.080368;  O=0.3938062999413425; P=0.3938062999413425; Q=104460.6053194707;  R=940145.4478752365;  S=0.3728236235608984;  T=0.3728236235608984  }
    8  = @{ E=3; F=1; G=30.939307;         H=92.81792100000001; I=0.03982438491082609;J=0.03982438491082609;K=3; L=1; M=117.044563;        N=351.133689;  O=0.3245365645427815; P=0.3245365645427815; Q=3621.277667337841;  R=32591.49900604057;  S=0.01292446906398888;  T=0.01292446906398888  }
    9  = @{ E=3; F=1; G=30.939307;         H=92.81792100000001; I=0.03982438491082609;J=0.03982438491082609;K=3; L=1; M=101.5800373333333; N=304.740112;  O=0.281657135515876;  P=0.281657135515876;  Q=3142.815960127462;  R=28285.34364114716;  S=0.01121682217766495;  T=0.01121682217766495  }
    10 = @{ E=3; F=1; G=30.939307;         H=92.81792100000001; I=0.03982438491082609;J=0.03982438491082609;K=3; L=1; M=142.0267893333333; N=426.080368;  O=0.3938062999413425; P=0.3938062999413425; Q=4394.210437408326;  R=39547.89393667493;  S=0.01568309366917225;  T=0.01568309366917225  }
}

$columns = @('E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T')

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $columns) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
